$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace `old` with `new` but scoped to a single paragraph's Range
# so that ambiguous substrings (e.g. "com" inside "comprehended") elsewhere
# in the document are not touched.
# ---------------------------------------------------------------------------
function Replace-InParagraph($paraIndex, $old, $new) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range.Duplicate
    $found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND (para $paraIndex): $old"
    }
}

# ---------------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------------
Replace-InParagraph 1 "The Enigma of Untapped Human Neuro-Potential" `
    "A Journey into the World of Chemistry: Exploring the Elements and Their Reactions"

# ---------------------------------------------------------------------------
# 2. Author
# ---------------------------------------------------------------------------
Replace-InParagraph 2 "Jasmine Alva" "Professor Samuel Hughes"

# ---------------------------------------------------------------------------
# 3. Email paragraph: rebuild wholesale as 5 runs (samuel / . / hughes@highschool / . / edu)
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$rng3 = $d.Range($p3.Range.Start, $p3.Range.End)
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr32 = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr>'
$emailXml = "<w:p $ns><w:pPr><w:jc w:val=`"center`"/></w:pPr>" +
    "<w:r>$rPr32<w:t>samuel</w:t></w:r>" +
    "<w:r>$rPr32<w:t>.</w:t></w:r>" +
    "<w:r>$rPr32<w:t>hughes@highschool</w:t></w:r>" +
    "<w:r>$rPr32<w:t>.</w:t></w:r>" +
    "<w:r>$rPr32<w:t>edu</w:t></w:r></w:p>"
$rng3.InsertXML($emailXml)

# ---------------------------------------------------------------------------
# 4. Body paragraph (paragraph 5) - replace each run's text
# ---------------------------------------------------------------------------
Replace-InParagraph 5 "Delving into the intricate workings of the human mind, we encounter an uncharted realm brimming with untapped potential--a hidden symphony of neurons yet to be fully comprehended" `
    "Chemistry, the science of matter and its properties, invites us on an intriguing journey to understand the world around us"

Replace-InParagraph 5 " Like a vast cosmic map, the human brain holds the key to understanding the deepest recesses of human consciousness, behavior, and creativity" `
    " From the seemingly mundane objects in our daily lives to the awe-inspiring phenomena in the universe, chemistry provides the framework to unravel the secrets hidden within the elements and their interactions"

Replace-InParagraph 5 " As we probe the mysteries of neural pathways, we unlock the secrets of learning, memory, and decision-making, painting a vivid tapestry of human experience" `
    " In this captivating voyage, we will delve into the fundamental concepts of chemistry, unravelling the mysteries of matter at its most basic level"

Replace-InParagraph 5 "From the intricate dance of synapses to the symphony of neurotransmitters, the human brain is an orchestra of biological wonders" `
    "As we embark on this odyssey, we will explore the fascinating realm of chemical elements - the building blocks of all matter"

Replace-InParagraph 5 " Billions of neurons, connected by trillions of synapses, orchestrate a ceaseless symphony of electrochemical signals, weaving together the fabric of our thoughts, emotions, and actions" `
    " We will uncover their unique properties, unravel the patterns that govern their behavior, and delve into the dynamic forces that shape their interactions"

Replace-InParagraph 5 " This remarkable organ possesses an astonishing capacity to learn, adapt, and create, constantly reshaping its neural landscape in response to new experiences" `
    " Moreover, we will witness the captivating spectacle of chemical reactions, where substances undergo transformations, revealing new compounds with distinct characteristics"

Replace-InParagraph 5 "Our understanding of the human brain remains a work in progress, yet the discoveries made thus far have illuminated the extraordinary capabilities of this enigmatic organ" `
    "Our exploration will lead us to appreciate the profound impact chemistry has on our lives"

Replace-InParagraph 5 " We have witnessed the plasticity of the brain, its remarkable ability to reorganize itself after injury, and its capacity for exceptional learning and memory" `
    " From the food we consume to the medicines that heal us, from the materials that clothe us to the fuels that power our world, chemistry plays an essential role in shaping our existence"

Replace-InParagraph 5 " These insights have inspired new approaches to treating neurological disorders and paved the way for groundbreaking advancements in artificial intelligence" `
    " Furthermore, we will examine the intricate connections between chemistry and other disciplines, such as biology, physics, and engineering, highlighting the interdisciplinary nature of science"

# ---------------------------------------------------------------------------
# 5. "Summary" heading (paragraph 6) - unchanged text, nothing to do
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 6. Summary paragraph (paragraph 7) - replace each run's text
# ---------------------------------------------------------------------------
Replace-InParagraph 7 "The human brain, with its vast reserves of untapped potential, stands as a testament to the remarkable complexity and resilience of life" `
    "Our journey into the realm of chemistry has unveiled the captivating world of elements and their reactions"

Replace-InParagraph 7 " As we delve deeper into the mysteries of neurobiology, we uncover the intricately woven tapestry of neurons, synapses, and neurotransmitters that orchestrate the symphony of human consciousness" `
    " We explored the fundamental building blocks of matter, unravelled the patterns that govern their behavior, and witnessed the transformative power of chemical reactions"

Replace-InParagraph 7 " Each new discovery brings us closer to comprehending the enigmatic enigma of the human mind, providing hope for novel treatments and transformative technologies that will shape the future of humanity" `
    " Moreover, we recognized the profound impact chemistry has on our daily lives and its interdisciplinary connections with other scientific fields"

Replace-InParagraph 7 " The boundless potential of the human brain remains an enduring source of wonder and awe, beckoning us to explore the vast frontiers of human consciousness" `
    " The study of chemistry, with its intricate tapestry of concepts and applications, invites us to continue our exploration into the enchanting realm of atoms, molecules, and reactions, unlocking further secrets of the universe we inhabit"

# ---------------------------------------------------------------------------
# 7. Font rename: TimesNewToman -> Times New Roman, document-wide
# ---------------------------------------------------------------------------
$full = $d.Range(0, $d.Content.End)
$full.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------------
# 8. Append a new empty paragraph at the very end of the document body
# ---------------------------------------------------------------------------
$endPos = $d.Content.End - 1
$endRange = $d.Range($endPos, $endPos)
$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

Write-Output "All edits applied"
